$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product rows appended at the bottom of the sheet (rows 47-58)
$newRows = @(
    @("ANA HICKMANN", "AH90020-A01", "https://oculosshop.com.br/cdn/shop/files/56266A.jpg?v=1692213144"),
    @("ANA HICKMANN", "AH90020-D21", "https://oculosshop.com.br/cdn/shop/files/56267A.jpg?v=1692213105"),
    @("ANA HICKMANN", "HI10025-05A", "https://oculosshop.com.br/cdn/shop/files/56256_grande.jpg?v=1734534988"),
    @("ANA HICKMANN", "HI10047-04A", "https://cdn.sistemawbuy.com.br/arquivos/de0f9d459e6b0537fbed094b2cdd7dba/produtos/665799776dea5/hi10047-04a-6659ba3638f1c_mini.jpg"),
    @("ANA HICKMANN", "HI10047-09A", "https://imgs.extra.com.br/1570709409/1xg.jpg?imwidth=500"),
    @("ANA HICKMANN", "HI10048-04A", "https://cdn.sistemawbuy.com.br/arquivos/de0f9d459e6b0537fbed094b2cdd7dba/produtos/668fde26dafaa/a-6728cdeaa10f1.jpg"),
    @("ANA HICKMANN", "HI10049-05A", "https://oculosshop.com.br/cdn/shop/files/101053_1.jpg?crop=center&height=1024&v=1727714814&width=1024"),
    @("ANA HICKMANN", "HI10050-09A", "https://loja.oticaademar.com.br/image/cache/catalog/Produtos/53633%20-%20Armacao%20Hickmann%20GARDENIA%201%20HI10050%2009A%2054-18%2001-800x800.jpg"),
    @("ANA HICKMANN", "HI10058-09C", "https://www.safira.com.br/cdn/imagens/produtos/det/armacao-feminina-ana-hickmann-hi10058-09c-4d3a6fa71cf53442c8574e716dd99fe8.jpg"),
    @("ANA HICKMANN", "HI10058-13C", "https://oculosshop.com.br/cdn/shop/files/101054-1.jpg?v=1727714989&width=2000"),
    @("ANA HICKMANN", "HI10060-01A", "https://oculosshop.com.br/cdn/shop/files/102933_1.jpg?v=1743794794"),
    @("ANA HICKMANN", "HI10060-06A", "https://http2.mlstatic.com/D_NQ_NP_983173-MLA95604340215_102025-O.webp")
)

$startRow = 47
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Turn the URL-valued cells into real hyperlinks (address == displayed text)
$ws.Hyperlinks.Add($ws.Range("C6"), $ws.Range("C6").Value)
$ws.Hyperlinks.Add($ws.Range("C10"), $ws.Range("C10").Value)
$ws.Hyperlinks.Add($ws.Range("C14"), $ws.Range("C14").Value)
$ws.Hyperlinks.Add($ws.Range("C16"), $ws.Range("C16").Value)
$ws.Hyperlinks.Add($ws.Range("C47"), $ws.Range("C47").Value)
